## Cobalt Users workbook update
## - Rename Sheet2 -> Emails, populate it with tr-anz tester emails
## - Append 28 new user rows (53-80) to the Users sheet, each with a
##   mailto hyperlink in column G and a left/right border on E:F
## - Trim Sheet3 back down to a single "Y" row
## - Touch up column widths / selections to match the authored file

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Users sheet (sheet1) - append the new accounts
# ---------------------------------------------------------------------------
$users = $wb.Worksheets.Item(1)

$newUsers = @(
    @{Row=53; User="SearchOpenWebUser1"; Email="SearchOpenWeb@mailinator.com "},
    @{Row=54; User="FFHUser1";           Email="FFHUser1@mailinator.com "},
    @{Row=55; User="FFHUser2";           Email="FFHUser2@mailinator.com"},
    @{Row=56; User="FFHUser3";           Email="FFHUser3@mailinator.com"},
    @{Row=57; User="FFHUser4";           Email="FFHUser4@mailinator.com"},
    @{Row=58; User="FrontEndUser1";      Email="FrontEndUser1@mailinator.com"},
    @{Row=59; User="FrontEndUser2";      Email="FrontEndUser2@mailinator.com"},
    @{Row=60; User="FrontEndUser3";      Email="FrontEndUser3@mailinator.com"},
    @{Row=61; User="FrontEndUser4";      Email="FrontEndUser4@mailinator.com"},
    @{Row=62; User="FrontEndUser5";      Email="FrontEndUser5@mailinator.com"},
    @{Row=63; User="FrontEndUser6";      Email="FrontEndUser6@mailinator.com"},
    @{Row=64; User="FrontEndUser7";      Email="FrontEndUser7@mailinator.com"},
    @{Row=65; User="FrontEndUser8";      Email="FrontEndUser8@mailinator.com"},
    @{Row=66; User="FrontEndUser9";      Email="FrontEndUser9@mailinator.com"},
    @{Row=67; User="FrontEndUser10";     Email="FrontEndUser10@mailinator.com"},
    @{Row=68; User="UrlUser1";           Email="UrlUser1@mailinator.com"},
    @{Row=69; User="UrlUser2";           Email="UrlUser2@mailinator.com"},
    @{Row=70; User="UrlUser3";           Email="UrlUser3@mailinator.com"},
    @{Row=71; User="LinkingUser1";       Email="LinkingUser1@mailinator.com "},
    @{Row=72; User="LoginUser1";         Email="LoginUser1@mailinator.com "},
    @{Row=73; User="LoginUser2";         Email="LoginUser2@mailinator.com "},
    @{Row=74; User="LoginUser3";         Email="LoginUser3@mailinator.com "},
    @{Row=75; User="LoginUser4";         Email="LoginUser4@mailinator.com "},
    @{Row=76; User="LoginUser5";         Email="LoginUser5@mailinator.com "},
    @{Row=77; User="LoginUser6";         Email="LoginUser6@mailinator.com "},
    @{Row=78; User="LoginUser7";         Email="LoginUser7@mailinator.com "},
    @{Row=79; User="CpetUser1";          Email="CpetUser1@mailinator.com "},
    @{Row=80; User="CpetUser2";          Email="CpetUser2@mailinator.com "}
)

foreach ($u in $newUsers) {
    $r = $u.Row
    $users.Cells.Item($r, 1).Value = $u.User           # A - UserName
    $users.Cells.Item($r, 2).Value = "Password1"        # B - Pasword
    $eCell = $users.Cells.Item($r, 5)
    $eCell.Value = "THIS IS IN USE 24/7 - DO NOT USE!"  # E - Description
    $fCell = $users.Cells.Item($r, 6)
    $fCell.Value = "N"                                  # F - Locked

    # Thin border on the left & right edges of the Description/Locked cells
    foreach ($cell in @($eCell, $fCell)) {
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(7).Weight = 2
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight = 2
    }

    $gCell = $users.Cells.Item($r, 7)
    $gCell.Value = $u.Email                             # G - Email
    [void]$users.Hyperlinks.Add($gCell, "mailto:" + $u.Email.Trim())
}

# Four trailing "empty but bordered" cells below the table (E81:E84)
foreach ($r in 81..84) {
    $cell = $users.Cells.Item($r, 5)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
}

# Column width tweaks on the Users sheet
$users.Columns.Item(1).ColumnWidth = 20.33   # A ~ 21.14
$users.Columns.Item(5).ColumnWidth = 38.5    # E ~ 39.29

[void]$users.Range("C80").Select()

# ---------------------------------------------------------------------------
# 2. Sheet2 -> Emails
# ---------------------------------------------------------------------------
$emails = $wb.Worksheets.Item(2)
$emails.Name = "Emails"

$emails.Cells.Item(1, 1).Value = "Email"
$emails.Cells.Item(1, 2).Value = "Password"
$emails.Cells.Item(2, 1).Value = "tr-anz-tester1@yandex.com"
$emails.Cells.Item(2, 2).Value = "tranztest"
$emails.Cells.Item(3, 1).Value = "tr-anz-tester2@yandex.com"
$emails.Cells.Item(3, 2).Value = "tranztest"

$emails.Columns.Item(1).ColumnWidth = 25.5   # ~ 26.29
$emails.Columns.Item(2).ColumnWidth = 13.0   # ~ 13.86

[void]$emails.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# 3. Sheet3 - trim back to a single "Y" row
# ---------------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Item(3)
[void]$sheet3.Range("A3:A7").EntireRow.Delete()

# Leave the Users sheet as the active sheet/selection, matching the source file
[void]$users.Activate()
[void]$users.Range("C80").Select()
